# Future Project Plans slide (slide 10): rewrite the bullet list in the
# content placeholder and remove the standalone "UPDATE" textbox.

$p = $ppt.ActivePresentation
$s = $p.Slides.Item(10)
$shp = $s.Shapes.Item(2)          # "Content Placeholder 2"
$tr = $shp.TextFrame.TextRange

$ldq = [char]8220                  # left double quotation mark
$rdq = [char]8221                  # right double quotation mark

# --- Paragraph 1: "Comparison by ... Turkiye, Indonesia, China ..." ---
# Rebuilt from three runs so the country name "Turkiye" stays its own run.
$para1 = $tr.Paragraphs(1, 1)
$para1.Text = "Comparison by " + $ldq + "Global South" + $rdq + " countries: Brazil, South Africa, "
$run2 = $para1.InsertAfter("Turkiye")
$run3 = $run2.InsertAfter(", Indonesia, China to see if the " + $ldq + "unaligned global south" + $rdq + " is less interested in divestment")

# --- Paragraph 2 (was "Duration of trending content") ---
$para2 = $tr.Paragraphs(2, 1)
$para2.Text = "Perhaps a scatter plot which shows the number of each letter grade by the 10 most represented countries in the study"

# --- Paragraph 3 (was "Identifying social media platforms in descriptions") ---
$para3 = $tr.Paragraphs(3, 1)
$para3.Text = "See if there is a difference between state-owner enterprises specifically, though this requires a more individuated look at each company"

# --- New paragraph 4, inserted after paragraph 3, before the trailing empty paragraph ---
$para4 = $para3.InsertAfter("`rCould group A & B and D & F to show greater aggregation")

# Remove the standalone "UPDATE" textbox shape entirely.
$s.Shapes.Item(3).Delete()
